$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 1174, shifting existing rows 1174..1259 down to 1175..1260
$ws.Rows.Item(1174).Insert(4)  # 4 = xlShiftDown

# Copy the date number format from the row below (now row 1175, originally row 1174)
$ws.Cells.Item(1174, 4).NumberFormat = $ws.Cells.Item(1175, 4).NumberFormat

# Populate the new row 1174 with the new data record
$ws.Cells.Item(1174, 1).Value = 8
$ws.Cells.Item(1174, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1174, 3).Value = "Coquimbo"
$ws.Cells.Item(1174, 4).Value = 45265
$ws.Cells.Item(1174, 5).Value = 4
$ws.Cells.Item(1174, 6).Value = 100112004
$ws.Cells.Item(1174, 7).Value = "Cebolla"
$ws.Cells.Item(1174, 8).Value = "Sin especificar"
$ws.Cells.Item(1174, 9).Value = "1a (cosecha)"
$ws.Cells.Item(1174, 10).Value = 2000
$ws.Cells.Item(1174, 11).Value = 13000
$ws.Cells.Item(1174, 12).Value = 14000
$ws.Cells.Item(1174, 13).Value = 13500
$ws.Cells.Item(1174, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(1174, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1174, 16).Value = 750
$ws.Cells.Item(1174, 17).Value = 18
$ws.Cells.Item(1174, 18).Value = "Hortaliza"
